$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 643, shifting existing rows 643:678 down to 644:679
$ws.Rows(643).Insert()

# Populate the newly inserted row 643 with the new "Florida King" record.
# Columns A,B,C,E,F,G,H,I,J keep the same values as the rest of this
# "Durazno" block (these were not part of the insert, so they need to be
# (re)written explicitly since Insert() leaves the new row blank).
$ws.Range("A643").Value = 6
$ws.Range("B643").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C643").Value = "Metropolitana"
$ws.Range("D643").Value = 44516
$ws.Range("E643").Value = 13
$ws.Range("F643").Value = "Fruta"
$ws.Range("G643").Value = 100103
$ws.Range("H643").Value = "Frutos de hueso (carozo)"
$ws.Range("I643").Value = 100103004
$ws.Range("J643").Value = "Durazno"
$ws.Range("K643").Value = "Florida King"
$ws.Range("L643").Value = "Primera"
$ws.Range("M643").Value = 6
$ws.Range("N643").Value = 450000
$ws.Range("O643").Value = 450000
$ws.Range("P643").Value = 450000
$ws.Range("Q643").Value = "`$/bins (420 kilos)"
$ws.Range("R643").Value = "Paine"
$ws.Range("S643").Value = 1071
$ws.Range("T643").Value = 420
